$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. date style on column A) from the last existing row (328)
# down across the newly appended rows (329-343), matching the style used throughout
# the sheet for the date column.
$ws.Range("A328:D328").Copy() | Out-Null
$ws.Range("A329:D343").PasteSpecial(-4122) | Out-Null

# Populate the new daily records (dates, new positives, 7-day rolling sum,
# and 7-day rolling sum per 100k inhabitants) through 9 August 2021.
$ws.Cells.Item(329, 1).Value = 44403
$ws.Cells.Item(329, 2).Value = 0
$ws.Cells.Item(329, 3).Value = 1
$ws.Cells.Item(329, 4).Value = 9.212344541685859

$ws.Cells.Item(330, 1).Value = 44404
$ws.Cells.Item(330, 2).Value = 0
$ws.Cells.Item(330, 3).Value = 1
$ws.Cells.Item(330, 4).Value = 9.212344541685859

$ws.Cells.Item(331, 1).Value = 44405
$ws.Cells.Item(331, 2).Value = 0
$ws.Cells.Item(331, 3).Value = 1
$ws.Cells.Item(331, 4).Value = 9.212344541685859

$ws.Cells.Item(332, 1).Value = 44406
$ws.Cells.Item(332, 2).Value = 0
$ws.Cells.Item(332, 3).Value = 1
$ws.Cells.Item(332, 4).Value = 9.212344541685859

$ws.Cells.Item(333, 1).Value = 44407
$ws.Cells.Item(333, 2).Value = 1
$ws.Cells.Item(333, 3).Value = 1
$ws.Cells.Item(333, 4).Value = 9.212344541685859

$ws.Cells.Item(334, 1).Value = 44408
$ws.Cells.Item(334, 2).Value = 0
$ws.Cells.Item(334, 3).Value = 1
$ws.Cells.Item(334, 4).Value = 9.212344541685859

$ws.Cells.Item(335, 1).Value = 44409
$ws.Cells.Item(335, 2).Value = 1
$ws.Cells.Item(335, 3).Value = 2
$ws.Cells.Item(335, 4).Value = 18.42468908337172

$ws.Cells.Item(336, 1).Value = 44410
$ws.Cells.Item(336, 2).Value = 0
$ws.Cells.Item(336, 3).Value = 2
$ws.Cells.Item(336, 4).Value = 18.42468908337172

$ws.Cells.Item(337, 1).Value = 44411
$ws.Cells.Item(337, 2).Value = 1
$ws.Cells.Item(337, 3).Value = 3
$ws.Cells.Item(337, 4).Value = 27.63703362505758

$ws.Cells.Item(338, 1).Value = 44412
$ws.Cells.Item(338, 2).Value = 0
$ws.Cells.Item(338, 3).Value = 3
$ws.Cells.Item(338, 4).Value = 27.63703362505758

$ws.Cells.Item(339, 1).Value = 44413
$ws.Cells.Item(339, 2).Value = 0
$ws.Cells.Item(339, 3).Value = 3
$ws.Cells.Item(339, 4).Value = 27.63703362505758

$ws.Cells.Item(340, 1).Value = 44414
$ws.Cells.Item(340, 2).Value = 1
$ws.Cells.Item(340, 3).Value = 3
$ws.Cells.Item(340, 4).Value = 27.63703362505758

$ws.Cells.Item(341, 1).Value = 44415
$ws.Cells.Item(341, 2).Value = 0
$ws.Cells.Item(341, 3).Value = 3
$ws.Cells.Item(341, 4).Value = 27.63703362505758

$ws.Cells.Item(342, 1).Value = 44416
$ws.Cells.Item(342, 2).Value = 1
$ws.Cells.Item(342, 3).Value = 3
$ws.Cells.Item(342, 4).Value = 27.63703362505758

$ws.Cells.Item(343, 1).Value = 44417
$ws.Cells.Item(343, 2).Value = 0
$ws.Cells.Item(343, 3).Value = 3
$ws.Cells.Item(343, 4).Value = 27.63703362505758
